$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 812.1875
$ws.Range("J53").Value = 427.6
$ws.Range("L53").Value = 427.6
$ws.Range("N53").Value = -1701.6

$ws.Range("H106").Value = 3825.7778
$ws.Range("I106").Value = 3825.7778
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3825.7778
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -3194.7778

$ws.Range("H112").Value = 3868.2812
$ws.Range("J112").Value = 4395.963
$ws.Range("L112").Value = 13187.889
$ws.Range("N112").Value = -15403.889

$ws.Range("H132").Value = 20792.766
$ws.Range("I132").Value = 21980.479
$ws.Range("J132").Value = 1789.3334
$ws.Range("K132").Value = 65941.43700000001
$ws.Range("L132").Value = 5368.0002
$ws.Range("M132").Value = -63411.43700000001
$ws.Range("N132").Value = -10428.0002

$ws.Range("H135").Value = 1667523.4
$ws.Range("I135").Value = 1667523.4
$ws.Range("K135").Value = 15007710.6
$ws.Range("M135").Value = -15005175.6

$ws.Range("H137").Value = 3420.1904
$ws.Range("I137").Value = 3542.7727
$ws.Range("J137").Value = 3285.35
$ws.Range("K137").Value = 10628.3181
$ws.Range("L137").Value = 9856.049999999999
$ws.Range("M137").Value = -8078.3181
$ws.Range("N137").Value = -14956.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2828848.5
$ws.Range("I32").Value = 3515075.2
$ws.Range("J32").Value = 34925.215
$ws.Range("K32").Value = 3515075.2
$ws.Range("L32").Value = 34925.215
$ws.Range("M32").Value = -3514788.2
$ws.Range("N32").Value = -35499.215

$ws.Range("H34").Value = 176749.75
$ws.Range("J34").Value = 68999
$ws.Range("L34").Value = 68999
$ws.Range("N34").Value = -69541

$ws.Range("H61").Value = 27780884
$ws.Range("I61").Value = 2202.5217
$ws.Range("K61").Value = 2202.5217
$ws.Range("M61").Value = -1990.5217

$ws.Range("H102").Value = 3593.8147
$ws.Range("I102").Value = 3043.0833
$ws.Range("K102").Value = 3043.0833
$ws.Range("M102").Value = -1421.0833

$ws.Range("H119").Value = 55938
$ws.Range("J119").Value = 55938
$ws.Range("L119").Value = 55938
$ws.Range("N119").Value = -65614

$ws.Range("H132").Value = 2631.2458
$ws.Range("I132").Value = 1094.238
$ws.Range("K132").Value = 3282.714
$ws.Range("M132").Value = -752.7139999999999

$ws.Range("H136").Value = 27780884
$ws.Range("I136").Value = 2202.5217
$ws.Range("K136").Value = 6607.5651
$ws.Range("M136").Value = -4057.5651

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 26316316
$ws.Range("I80").Value = 55556064
$ws.Range("J80").Value = 542.8
$ws.Range("K80").Value = 55556064
$ws.Range("L80").Value = 542.8
$ws.Range("M80").Value = -55555066
$ws.Range("N80").Value = -2538.8

$ws.Range("H83").Value = 26316316
$ws.Range("I83").Value = 55556064
$ws.Range("J83").Value = 542.8
$ws.Range("K83").Value = 277780320
$ws.Range("L83").Value = 2714
$ws.Range("M83").Value = -277775328
$ws.Range("N83").Value = -12698

$ws.Range("H113").Value = 4962.3335
$ws.Range("I113").Value = 4962.3335
$ws.Range("K113").Value = 4962.3335
$ws.Range("M113").Value = -2792.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5205.1226
$ws.Range("I31").Value = 1798.2354
$ws.Range("J31").Value = 7015.0312
$ws.Range("K31").Value = 1798.2354
$ws.Range("L31").Value = 7015.0312
$ws.Range("M31").Value = -1503.2354
$ws.Range("N31").Value = -7605.0312

$ws.Range("H34").Value = 5205.1226
$ws.Range("I34").Value = 1798.2354
$ws.Range("J34").Value = 7015.0312
$ws.Range("K34").Value = 1798.2354
$ws.Range("L34").Value = 7015.0312
$ws.Range("M34").Value = -1596.2354
$ws.Range("N34").Value = -7419.0312

$ws.Range("H132").Value = 3984.9302
$ws.Range("I132").Value = 2683.8928
$ws.Range("J132").Value = 6413.533
$ws.Range("K132").Value = 8051.678400000001
$ws.Range("L132").Value = 19240.599
$ws.Range("M132").Value = -5521.678400000001
$ws.Range("N132").Value = -24300.599

$ws.Range("H141").Value = 64516.4
$ws.Range("J141").Value = 64516.4
$ws.Range("L141").Value = 64516.4
$ws.Range("N141").Value = -74876.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1232
$ws.Range("I5").Value = 783.1667
$ws.Range("K5").Value = 2349.5001
$ws.Range("M5").Value = -2237.5001

$ws.Range("H39").Value = 13399.143
$ws.Range("J39").Value = 13399.143
$ws.Range("L39").Value = 40197.429
$ws.Range("N39").Value = -40785.429

$ws.Range("H76").Value = 250002700
$ws.Range("I76").Value = 250002700
$ws.Range("K76").Value = 750008100
$ws.Range("M76").Value = -750007717

$ws.Range("H79").Value = 250002700
$ws.Range("I79").Value = 250002700
$ws.Range("K79").Value = 750008100
$ws.Range("M79").Value = -750006774

$ws.Range("H86").Value = 398.5
$ws.Range("J86").Value = 397
$ws.Range("L86").Value = 1191
$ws.Range("N86").Value = -3563

$ws.Range("H89").Value = 398.5
$ws.Range("J89").Value = 397
$ws.Range("L89").Value = 3573
$ws.Range("N89").Value = -15429

$ws.Range("H122").Value = 4716021
$ws.Range("I122").Value = 4716021
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 42444189
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -42441739

$ws.Range("H135").Value = 1232
$ws.Range("I135").Value = 783.1667
$ws.Range("K135").Value = 7048.5003
$ws.Range("M135").Value = -4513.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 42962
$ws.Range("J33").Value = 42962
$ws.Range("L33").Value = 42962
$ws.Range("N33").Value = -43466

$ws.Range("H126").Value = 4802
$ws.Range("I126").Value = 3504.8333
$ws.Range("J126").Value = 5774.875
$ws.Range("K126").Value = 10514.4999
$ws.Range("L126").Value = 17324.625
$ws.Range("M126").Value = -8044.499899999999
$ws.Range("N126").Value = -22264.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4789.4443
$ws.Range("J7").Value = 5597.3
$ws.Range("L7").Value = 5597.3
$ws.Range("N7").Value = -5821.3

$ws.Range("H22").Value = 1347.1765
$ws.Range("I22").Value = 442.57144
$ws.Range("J22").Value = 1980.4
$ws.Range("K22").Value = 442.57144
$ws.Range("L22").Value = 1980.4
$ws.Range("M22").Value = -147.57144
$ws.Range("N22").Value = -2570.4

$ws.Range("H27").Value = 1347.1765
$ws.Range("I27").Value = 442.57144
$ws.Range("J27").Value = 1980.4
$ws.Range("K27").Value = 442.57144
$ws.Range("L27").Value = 1980.4
$ws.Range("M27").Value = -335.57144
$ws.Range("N27").Value = -2194.4

$ws.Range("H46").Value = 1047265
$ws.Range("I46").Value = 2030116
$ws.Range("J46").Value = 2985.9375
$ws.Range("K46").Value = 2030116
$ws.Range("L46").Value = 2985.9375
$ws.Range("M46").Value = -2029928
$ws.Range("N46").Value = -3361.9375

$ws.Range("H82").Value = 1006985.75
$ws.Range("I82").Value = 1677016
$ws.Range("J82").Value = 1940.3334
$ws.Range("K82").Value = 1677016
$ws.Range("L82").Value = 1940.3334
$ws.Range("M82").Value = -1676655
$ws.Range("N82").Value = -2662.3334

$ws.Range("H85").Value = 1006985.75
$ws.Range("I85").Value = 1677016
$ws.Range("J85").Value = 1940.3334
$ws.Range("K85").Value = 1677016
$ws.Range("L85").Value = 1940.3334
$ws.Range("M85").Value = -1675768
$ws.Range("N85").Value = -4436.3334

$ws.Range("H107").Value = 4186.75
$ws.Range("I107").Value = 4186.75
$ws.Range("K107").Value = 4186.75
$ws.Range("M107").Value = -2266.75

$ws.Range("H126").Value = 4789.4443
$ws.Range("J126").Value = 5597.3
$ws.Range("L126").Value = 16791.9
$ws.Range("N126").Value = -21731.9

$ws.Range("H132").Value = 8337422
$ws.Range("I132").Value = 15154428
$ws.Range("J132").Value = 5526.1113
$ws.Range("K132").Value = 45463284
$ws.Range("L132").Value = 16578.3339
$ws.Range("M132").Value = -45460754
$ws.Range("N132").Value = -21638.3339

$ws.Range("H136").Value = 7217.5605
$ws.Range("I136").Value = 2437.025
$ws.Range("J136").Value = 14572.23
$ws.Range("K136").Value = 7311.075000000001
$ws.Range("L136").Value = 43716.69
$ws.Range("M136").Value = -4761.075000000001
$ws.Range("N136").Value = -48816.69

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 12346377
$ws.Range("I107").Value = 579.0625
$ws.Range("K107").Value = 1737.1875
$ws.Range("M107").Value = 182.8125
